$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings are not
# auto-converted to numbers by Excel, then restore the default style
# so no spurious formatting differences are introduced.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = '29.501.11'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.906.03'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '325.74'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").Value = '0.4843'
$ws.Range("E7").Value = '  +3.71%  '
$ws.Range("D8").Value = '0.4072'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '0.08130'
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '1.010'
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").Value = '23.50'
$ws.Range("E11").Value = '  +5.48%  '
$ws.Range("D12").Value = '1.902.50'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").Value = '6.023'
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '7.093'
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '90.45'
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '0.06759'
$ws.Range("E17").Value = '  +2.50%  '
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").Value = '17.69'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '29.509.10'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '5.572'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").Value = '11.79'
$ws.Range("E23").Value = '  +2.82%  '
$ws.Range("D24").Value = '2.164'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("D25").Value = '2.119.77'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D26").Value = '154.03'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '20.07'
$ws.Range("E27").Value = '  +1.72%  '
$ws.Range("D28").Value = '6.227'
$ws.Range("E28").Value = '  +9.20%  '
$ws.Range("D29").Value = '2.103'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").Value = '119.06'
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").Value = '1.036'
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("D32").Value = '0.09552'
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = '5.525'
$ws.Range("E33").Value = '  +2.57%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.393'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.550'
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("D36").Value = '0.02266'
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").Value = '0.06105'
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").Value = '0.5940'
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = '7.914'
$ws.Range("E40").Value = '  -5.52%  '
$ws.Range("D41").Value = '10.36'
$ws.Range("E41").Value = '  +2.72%  '
$ws.Range("D42").Value = '0.1857'
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.410'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = '1.283'
$ws.Range("E44").Value = '  -1.57%  '
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").Value = '12.44'
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("D47").Value = '0.5559'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").Value = '1.957'
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("D49").Value = '115.11'
$ws.Range("E49").Value = '  +1.64%  '
$ws.Range("D50").Value = '72.65'
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").Value = '1.056'
$ws.Range("E51").Value = '  +2.69%  '

$colD.Style = "Normal"
